$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written/kept as text, even when it
# looks numeric (e.g. "47.53"), so Excel does not silently convert
# it into a floating point number. The cell style is restored to
# "Normal" afterwards so no stray style index lingers on the cell.
function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# --- Apply the updated crypto price/volume figures ---
Set-TextValue $ws.Range('D2') '43.729.66'
$ws.Range('E2').Value = '  +0.42%  '
Set-TextValue $ws.Range('D3') '2.285.68'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.52%  '
Set-TextValue $ws.Range('D5') '110.17'
$ws.Range('E5').Value = '  +15.27%  '
Set-TextValue $ws.Range('D6') '267.66'
$ws.Range('E6').Value = '  -0.02%  '
Set-TextValue $ws.Range('D7') '0.624'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +0.33%  '
Set-TextValue $ws.Range('D9') '0.614'
$ws.Range('E9').Value = '  +0.99%  '
Set-TextValue $ws.Range('D10') '47.53'
Set-TextValue $ws.Range('D11') '0.0945'
$ws.Range('E11').Value = '  +1.51%  '
Set-TextValue $ws.Range('D12') '9.01'
$ws.Range('E12').Value = '  +13.97%  '
$ws.Range('E13').Value = '  +1.07%  '
Set-TextValue $ws.Range('D14') '15.73'
$ws.Range('E14').Value = '  +3.15%  '
Set-TextValue $ws.Range('D15') '2.630.03'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('E16').Value = '  -0.16%  '
Set-TextValue $ws.Range('D17') '2.287.39'
$ws.Range('E17').Value = '  +0.23%  '
Set-TextValue $ws.Range('D18') '43.599.78'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  +0.19%  '
Set-TextValue $ws.Range('D20') '6.80'
$ws.Range('E20').Value = '  +9.50%  '
Set-TextValue $ws.Range('D21') '72.12'
$ws.Range('E21').Value = '  +0.18%  '
Set-TextValue $ws.Range('D22') '2.46'
$ws.Range('E22').Value = '  -5.17%  '
Set-TextValue $ws.Range('D23') '232.16'
$ws.Range('E23').Value = '  +0.00%  '
Set-TextValue $ws.Range('D24') '9.79'
$ws.Range('E24').Value = '  +7.68%  '
Set-TextValue $ws.Range('D25') '2.76'
$ws.Range('E25').Value = '  +8.81%  '
$ws.Range('E26').Value = '  -0.04%  '
Set-TextValue $ws.Range('D27') '11.69'
$ws.Range('E27').Value = '  +4.69%  '
Set-TextValue $ws.Range('D28') '3.94'
$ws.Range('E28').Value = '  +1.12%  '
Set-TextValue $ws.Range('D29') '41.74'
$ws.Range('E29').Value = '  +4.44%  '
Set-TextValue $ws.Range('D30') '3.39'
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('E31').Value = '  -0.71%  '
Set-TextValue $ws.Range('D32') '175.74'
$ws.Range('E32').Value = '  +0.43%  '
Set-TextValue $ws.Range('D33') '21.54'
$ws.Range('E33').Value = '  -0.92%  '
Set-TextValue $ws.Range('D34') '0.0925'
$ws.Range('E35').Value = '  +4.98%  '
$ws.Range('E36').Value = '  +1.17%  '
Set-TextValue $ws.Range('D37') '4.71'
$ws.Range('E37').Value = '  +8.14%  '
$ws.Range('E38').Value = '  +4.02%  '
Set-TextValue $ws.Range('D39') '0.107'
$ws.Range('E39').Value = '  +0.36%  '
Set-TextValue $ws.Range('D40') '3.79'
$ws.Range('E40').Value = '  +13.39%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D41') '13.78'
$ws.Range('E41').Value = '  +12.25%  '
$ws.Range('E42').Value = '  +1.75%  '
$ws.Range('B43').Value = 'LidoDAOToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D43') '2.40'
$ws.Range('E43').Value = '  +4.31%  '
Set-TextValue $ws.Range('D44') '73.01'
$ws.Range('E44').Value = '  +11.07%  '
Set-TextValue $ws.Range('D45') '6.18'
$ws.Range('E45').Value = '  +19.72%  '
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('E47').Value = '  +2.06%  '
Set-TextValue $ws.Range('D48') '8.82'
$ws.Range('E48').Value = '  +0.84%  '
Set-TextValue $ws.Range('D49') '102.36'
$ws.Range('E49').Value = '  +6.35%  '
Set-TextValue $ws.Range('D50') '0.0989'
$ws.Range('E50').Value = '  -2.62%  '
$ws.Range('E51').Value = '  +2.82%  '
